$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 6: add date, proveedor (GUSTAVO), importe vendido, fecha de pago, importe cobrado
$ws.Range("A6").Value = 44456
$ws.Range("D6").Value = "GUSTAVO"
$ws.Range("E6").Value = 2642
$ws.Range("F6").Value = 44457
$ws.Range("G6").Value = 2642

# Row 7: add date, proveedor (GUSTAVO), importe vendido (no fecha de pago / importe cobrado yet)
$ws.Range("A7").Value = 44457
$ws.Range("D7").Value = "GUSTAVO"
$ws.Range("E7").Value = 8110

# Update the selection to match the diff (active cell F7)
$ws.Range("F7").Select()

$wb.Save()
